$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the Treatment query text in B5: remove the redundant CONCAT(...) wrapper ---
$treatmentCell = $ws.Cells.Item(5, 2)
$oldText = $treatmentCell.Value()
$newText = $oldText.Replace("CONCAT(REPLACE(trt.treatment_agent, ';', ', '))", "REPLACE(trt.treatment_agent, ';', ', ')")
$treatmentCell.Value = $newText

# Re-touch the font so the cell picks up a refreshed style record (mirrors the
# font/style churn produced when the author re-typed this cell's contents).
$treatmentCell.Font.ThemeColor = 1
$treatmentCell.Font.Size = 12
$treatmentCell.WrapText = $true

# --- Update the saved view/selection state: select B2 (clears any frozen/scrolled topLeftCell) ---
$ws.Activate()
$ws.Range("B2").Select()

$wb.Save()
